# Update the frequency table (rows 2-5, columns B-X) with newly computed
# proportions for the re-run used for publication.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.00436363636363636, 0.00218181818181818, 0.378181818181818, 0.0269090909090909, 0.00290909090909091, 0.985454545454545, 0.008, 0.00727272727272727, 0.00363636363636364, 0.356363636363636, 0.968727272727273, 0.00436363636363636, 0.0327272727272727, 0.00509090909090909, 0.979636363636364, 0.0385454545454545, 0.00436363636363636, 0.979636363636364, 0.000727272727272727, 0.0443636363636364, 0.0152727272727273, 0.277090909090909, 0.00945454545454546)
    3 = @(0.965818181818182, 0.976, 0.00290909090909091, 0.00363636363636364, 0.992727272727273, 0.0123636363636364, 0.00872727272727273, 0, 0, 0.000727272727272727, 0.0101818181818182, 0.00436363636363636, 0.00509090909090909, 0.00436363636363636, 0.00727272727272727, 0.0625454545454545, 0.00727272727272727, 0.00145454545454545, 0.973090909090909, 0.850181818181818, 0.962909090909091, 0.00145454545454545, 0.00218181818181818)
    4 = @(0.00945454545454546, 0.008, 0.0123636363636364, 0.968, 0.00290909090909091, 0.00218181818181818, 0.00581818181818182, 0.992, 0.994909090909091, 0.64, 0.0196363636363636, 0.0196363636363636, 0.928727272727273, 0.984, 0.00727272727272727, 0.330909090909091, 0.977454545454545, 0.0123636363636364, 0.0167272727272727, 0.0203636363636364, 0.0101818181818182, 0.526545454545455, 0.985454545454545)
    5 = @(0.0203636363636364, 0.0138181818181818, 0.606545454545454, 0.000727272727272727, 0.00145454545454545, 0, 0.977454545454545, 0.000727272727272727, 0.00145454545454545, 0.000727272727272727, 0.00145454545454545, 0.971636363636364, 0.032, 0.00654545454545455, 0.00581818181818182, 0.568, 0.0109090909090909, 0.00654545454545455, 0.00945454545454546, 0.0850909090909091, 0.0116363636363636, 0.193454545454545, 0.00218181818181818)
}

# Columns B (2) through X (24)
foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
